$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.856.03'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.907.76'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.10'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5020'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.52%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07281'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9109'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.21%  '
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07646'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').Value = '1.877.60'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.505'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.89'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008734'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '27.882.22'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.181'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.81'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.594'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.882'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.219'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.55%  '
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.50'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.929'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09037'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.207'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.228'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.776'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.20%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02083'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.526'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.093'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5537'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.020'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05275'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.894'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.482'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1516'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '111.74'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4822'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.633'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.55'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06058'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9075'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.93%  '
